# Update header row column names: remove spaces and accented characters.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "Nombre_comun"
$ws.Range("F1").Value = "Genero"
$ws.Range("G1").Value = "Epiteto_especifico"
$ws.Range("B1").Value = "Filo_o_division"

# Update the active selection on the sheet (matches the saved view state in the diff).
$ws.Range("I21").Select()

$wb.Save()
